$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to remain text storage (matches original inlineStr text cells)
# so that numeric-looking values like "0.7145" are not auto-converted to numbers.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.177.41"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.861.87"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "0.7145"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "240.23"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "0.07735"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "0.3068"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "24.86"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").Value = "0.08249"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.842.85"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").Value = "0.7148"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "5.205"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "90.16"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "29.181.57"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "5.840"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "243.02"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "0.000007781"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "2.118.58"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "13.11"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "7.916"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "0.1591"
$ws.Range("E25").Value = "  +8.93%  "
$ws.Range("D26").Value = "162.15"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "8.890"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "18.16"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "1.331"
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("D30").Value = "1.493"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "4.342"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").Value = "4.083"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").Value = "0.05181"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "1.912"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "1.172"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").Value = "0.7272"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("D37").Value = "2.678"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "2.695"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "1.155.74"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("D41").Value = "0.9012"
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("D42").Value = "6.092"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("D43").Value = "72.00"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "101.45"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D48").Value = "1.760"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").Value = "9.273"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").Value = "2.867"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").Value = "  -0.55%  "

# Row 46/47 swap: RocketPoolETH <-> Mantle, with updated figures
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.5282"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.012.53"
$ws.Range("E47").Value = "  -0.49%  "

# Restore default (no explicit) style so the saved XML has no "s" attribute,
# matching the original cell styling which was untouched by this data update.
$priceVolRange.Style = "Normal"
